$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845229
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45164.45833333334
$arr[0,3] = "Ayr"
$arr[0,4] = "Dundee Utd"
$arr[0,5] = 0
$arr[0,6] = 3
$arr[0,7] = "A"
$arr[0,8] = 2.3
$arr[0,9] = 3.3
$arr[0,10] = 2.8
$arr[0,11] = 3.1
$arr[0,12] = 3.4
$arr[0,13] = 2.25
$arr[0,14] = 0.25
$arr[0,15] = 1.875
$arr[0,16] = 1.925
$arr[0,17] = 2.25
$arr[0,18] = 1.8
$arr[0,19] = 2
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 1.25
$arr[0,23] = -1
$arr[0,24] = 0.925
$arr[0,25] = 0.8
$arr[0,26] = -1
$ws.Range("B15:AB15").Value2 = $arr

# Row 17
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845230
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45164.45833333334
$arr[0,3] = "Dunfermline"
$arr[0,4] = "Raith"
$arr[0,5] = 0
$arr[0,6] = 1
$arr[0,7] = "A"
$arr[0,8] = 2.1
$arr[0,9] = 3.3
$arr[0,10] = 3.2
$arr[0,11] = 2.45
$arr[0,12] = 3.4
$arr[0,13] = 2.8
$arr[0,14] = 0
$arr[0,15] = 1.8
$arr[0,16] = 2.05
$arr[0,17] = 2.5
$arr[0,18] = 2.05
$arr[0,19] = 1.8
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 1.8
$arr[0,23] = -1
$arr[0,24] = 1.05
$arr[0,25] = -1
$arr[0,26] = 0.8
$ws.Range("B17:AB17").Value2 = $arr

# Row 18
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845237
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45171.45833333334
$arr[0,3] = "Raith"
$arr[0,4] = "Queens Park"
$arr[0,5] = 3
$arr[0,6] = 2
$arr[0,7] = "H"
$arr[0,8] = 2.375
$arr[0,9] = 3.5
$arr[0,10] = 2.6
$arr[0,11] = 2.1
$arr[0,12] = 3.6
$arr[0,13] = 3.1
$arr[0,14] = -0.25
$arr[0,15] = 1.875
$arr[0,16] = 1.925
$arr[0,17] = 2.75
$arr[0,18] = 1.975
$arr[0,19] = 1.825
$arr[0,20] = 1.1
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.875
$arr[0,24] = -1
$arr[0,25] = 0.9750000000000001
$arr[0,26] = -1
$ws.Range("B18:AB18").Value2 = $arr

# Row 19
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845236
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45171.45833333334
$arr[0,3] = "Morton"
$arr[0,4] = "Partick"
$arr[0,5] = 1
$arr[0,6] = 4
$arr[0,7] = "A"
$arr[0,8] = 2.3
$arr[0,9] = 3.6
$arr[0,10] = 2.6
$arr[0,11] = 2.7
$arr[0,12] = 3.75
$arr[0,13] = 2.3
$arr[0,14] = 0.25
$arr[0,15] = 1.775
$arr[0,16] = 2.025
$arr[0,17] = 2.5
$arr[0,18] = 1.8
$arr[0,19] = 2
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 1.3
$arr[0,23] = -1
$arr[0,24] = 1.025
$arr[0,25] = 0.8
$arr[0,26] = -1
$ws.Range("B19:AB19").Value2 = $arr

# Row 24
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845242
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45185.45833333334
$arr[0,3] = "Raith"
$arr[0,4] = "Inverness CT"
$arr[0,5] = 1
$arr[0,6] = 0
$arr[0,7] = "H"
$arr[0,8] = 2.1
$arr[0,9] = 3.5
$arr[0,10] = 2.8
$arr[0,11] = 1.8
$arr[0,12] = 3.75
$arr[0,13] = 4
$arr[0,14] = -0.5
$arr[0,15] = 1.8
$arr[0,16] = 2
$arr[0,17] = 2.5
$arr[0,18] = 1.875
$arr[0,19] = 1.925
$arr[0,20] = 0.8
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.8
$arr[0,24] = -1
$arr[0,25] = -1
$arr[0,26] = 0.925
$ws.Range("B24:AB24").Value2 = $arr

# Row 25
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845241
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45185.45833333334
$arr[0,3] = "Queens Park"
$arr[0,4] = "Dunfermline"
$arr[0,5] = 0
$arr[0,6] = 2
$arr[0,7] = "A"
$arr[0,8] = 2.2
$arr[0,9] = 3.4
$arr[0,10] = 2.75
$arr[0,11] = 2.55
$arr[0,12] = 3.4
$arr[0,13] = 2.5
$arr[0,14] = 0
$arr[0,15] = 1.925
$arr[0,16] = 1.875
$arr[0,17] = 2.5
$arr[0,18] = 1.85
$arr[0,19] = 1.95
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 1.5
$arr[0,23] = -1
$arr[0,24] = 0.875
$arr[0,25] = -1
$arr[0,26] = 0.95
$ws.Range("B25:AB25").Value2 = $arr

# Row 26
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845240
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45185.45833333334
$arr[0,3] = "Dundee Utd"
$arr[0,4] = "Morton"
$arr[0,5] = 1
$arr[0,6] = 1
$arr[0,7] = "D"
$arr[0,8] = 1.571
$arr[0,9] = 4
$arr[0,10] = 4.333
$arr[0,11] = 1.55
$arr[0,12] = 4
$arr[0,13] = 5.25
$arr[0,14] = -1
$arr[0,15] = 1.95
$arr[0,16] = 1.85
$arr[0,17] = 2.75
$arr[0,18] = 1.95
$arr[0,19] = 1.85
$arr[0,20] = -1
$arr[0,21] = 3
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.8500000000000001
$arr[0,25] = -1
$arr[0,26] = 0.8500000000000001
$ws.Range("B26:AB26").Value2 = $arr

# Row 27
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845238
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45185.45833333334
$arr[0,3] = "Arbroath"
$arr[0,4] = "Airdrieonians"
$arr[0,5] = 4
$arr[0,6] = 0
$arr[0,7] = "H"
$arr[0,8] = 2.25
$arr[0,9] = 3.25
$arr[0,10] = 2.75
$arr[0,11] = 2.625
$arr[0,12] = 3.2
$arr[0,13] = 2.625
$arr[0,14] = 0
$arr[0,15] = 1.925
$arr[0,16] = 1.875
$arr[0,17] = 2.5
$arr[0,18] = 1.875
$arr[0,19] = 1.925
$arr[0,20] = 1.625
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.925
$arr[0,24] = -1
$arr[0,25] = 0.875
$arr[0,26] = -1
$ws.Range("B27:AB27").Value2 = $arr

# Row 74
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845295
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45276.5
$arr[0,3] = "Ayr"
$arr[0,4] = "Airdrieonians"
$arr[0,5] = 1
$arr[0,6] = 0
$arr[0,7] = "H"
$arr[0,8] = 2.1
$arr[0,9] = 3.5
$arr[0,10] = 3
$arr[0,11] = 2.2
$arr[0,12] = 3.2
$arr[0,13] = 3.1
$arr[0,14] = -0.25
$arr[0,15] = 2
$arr[0,16] = 1.85
$arr[0,17] = 2.25
$arr[0,18] = 2.025
$arr[0,19] = 1.825
$arr[0,20] = 1.2
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 1
$arr[0,24] = -1
$arr[0,25] = -1
$arr[0,26] = 0.825
$ws.Range("B74:AB74").Value2 = $arr

# Row 76
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845296
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45276.5
$arr[0,3] = "Dundee Utd"
$arr[0,4] = "Raith"
$arr[0,5] = 0
$arr[0,6] = 1
$arr[0,7] = "A"
$arr[0,8] = 1.8
$arr[0,9] = 3.4
$arr[0,10] = 4
$arr[0,11] = 1.7
$arr[0,12] = 3.4
$arr[0,13] = 4.5
$arr[0,14] = -0.75
$arr[0,15] = 2
$arr[0,16] = 1.8
$arr[0,17] = 2.5
$arr[0,18] = 2
$arr[0,19] = 1.8
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 3.5
$arr[0,23] = -1
$arr[0,24] = 0.8
$arr[0,25] = -1
$arr[0,26] = 0.8
$ws.Range("B76:AB76").Value2 = $arr

# Row 78
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845298
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45276.5
$arr[0,3] = "Inverness CT"
$arr[0,4] = "Arbroath"
$arr[0,5] = 1
$arr[0,6] = 2
$arr[0,7] = "A"
$arr[0,8] = 1.75
$arr[0,9] = 3.6
$arr[0,10] = 4
$arr[0,11] = 1.571
$arr[0,12] = 3.8
$arr[0,13] = 5
$arr[0,14] = -1
$arr[0,15] = 2.025
$arr[0,16] = 1.775
$arr[0,17] = 2.5
$arr[0,18] = 1.825
$arr[0,19] = 1.975
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 4
$arr[0,23] = -1
$arr[0,24] = 0.7749999999999999
$arr[0,25] = 0.825
$arr[0,26] = -1
$ws.Range("B78:AB78").Value2 = $arr

# Row 79
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845304
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45282.69791666666
$arr[0,3] = "Raith"
$arr[0,4] = "Ayr"
$arr[0,5] = 4
$arr[0,6] = 4
$arr[0,7] = "D"
$arr[0,8] = 1.666
$arr[0,9] = 4
$arr[0,10] = 4.75
$arr[0,11] = 1.533
$arr[0,12] = 4.2
$arr[0,13] = 5.5
$arr[0,14] = -1
$arr[0,15] = 1.925
$arr[0,16] = 1.875
$arr[0,17] = 2.75
$arr[0,18] = 1.975
$arr[0,19] = 1.825
$arr[0,20] = -1
$arr[0,21] = 3.2
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.875
$arr[0,25] = 0.9750000000000001
$arr[0,26] = -1
$ws.Range("B79:AB79").Value2 = $arr

# Row 80
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845303
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45282.69791666666
$arr[0,3] = "Queens Park"
$arr[0,4] = "Dundee Utd"
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = "D"
$arr[0,8] = 6
$arr[0,9] = 4.333
$arr[0,10] = 1.5
$arr[0,11] = 6
$arr[0,12] = 4.333
$arr[0,13] = 1.5
$arr[0,14] = 1.25
$arr[0,15] = 1.775
$arr[0,16] = 2.025
$arr[0,17] = 3
$arr[0,18] = 2
$arr[0,19] = 1.8
$arr[0,20] = -1
$arr[0,21] = 3.333
$arr[0,22] = -1
$arr[0,23] = 0.7749999999999999
$arr[0,24] = -1
$arr[0,25] = -1
$arr[0,26] = 0.8
$ws.Range("B80:AB80").Value2 = $arr

# Row 85
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845309
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45290.5
$arr[0,3] = "Raith"
$arr[0,4] = "Arbroath"
$arr[0,5] = 2
$arr[0,6] = 2
$arr[0,7] = "D"
$arr[0,8] = 1.444
$arr[0,9] = 4.2
$arr[0,10] = 5.5
$arr[0,11] = 1.4
$arr[0,12] = 4.333
$arr[0,13] = 6
$arr[0,14] = -1.25
$arr[0,15] = 1.875
$arr[0,16] = 1.925
$arr[0,17] = 3
$arr[0,18] = 1.9
$arr[0,19] = 1.9
$arr[0,20] = -1
$arr[0,21] = 3.333
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.925
$arr[0,25] = 0.8999999999999999
$arr[0,26] = -1
$ws.Range("B85:AB85").Value2 = $arr

# Row 86
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845307
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45290.5
$arr[0,3] = "Inverness CT"
$arr[0,4] = "Morton"
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = "D"
$arr[0,8] = 2.2
$arr[0,9] = 3.4
$arr[0,10] = 2.75
$arr[0,11] = 2.1
$arr[0,12] = 3.1
$arr[0,13] = 3.1
$arr[0,14] = -0.25
$arr[0,15] = 2.025
$arr[0,16] = 1.825
$arr[0,17] = 2
$arr[0,18] = 1.95
$arr[0,19] = 1.9
$arr[0,20] = -1
$arr[0,21] = 2.1
$arr[0,22] = -1
$arr[0,23] = -0.5
$arr[0,24] = 0.4125
$arr[0,25] = -1
$arr[0,26] = 0.8999999999999999
$ws.Range("B86:AB86").Value2 = $arr

# Row 87
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845305
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45290.5
$arr[0,3] = "Ayr"
$arr[0,4] = "Dunfermline"
$arr[0,5] = 2
$arr[0,6] = 2
$arr[0,7] = "D"
$arr[0,8] = 2.1
$arr[0,9] = 3.4
$arr[0,10] = 2.9
$arr[0,11] = 2.375
$arr[0,12] = 3.3
$arr[0,13] = 2.55
$arr[0,14] = 0
$arr[0,15] = 1.85
$arr[0,16] = 1.95
$arr[0,17] = 2.5
$arr[0,18] = 1.95
$arr[0,19] = 1.85
$arr[0,20] = -1
$arr[0,21] = 2.3
$arr[0,22] = -1
$arr[0,23] = 0
$arr[0,24] = 0
$arr[0,25] = 0.95
$arr[0,26] = -1
$ws.Range("B87:AB87").Value2 = $arr

# Row 88
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845308
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45290.5
$arr[0,3] = "Queens Park"
$arr[0,4] = "Airdrieonians"
$arr[0,5] = 1
$arr[0,6] = 2
$arr[0,7] = "A"
$arr[0,8] = 2.45
$arr[0,9] = 3.4
$arr[0,10] = 2.375
$arr[0,11] = 2.4
$arr[0,12] = 3.4
$arr[0,13] = 2.4
$arr[0,14] = 0
$arr[0,15] = 1.9
$arr[0,16] = 1.9
$arr[0,17] = 2.5
$arr[0,18] = 1.875
$arr[0,19] = 1.925
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 1.4
$arr[0,23] = -1
$arr[0,24] = 0.8999999999999999
$arr[0,25] = 0.875
$arr[0,26] = -1
$ws.Range("B88:AB88").Value2 = $arr

# Row 89
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845312
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45293.5
$arr[0,3] = "Morton"
$arr[0,4] = "Ayr"
$arr[0,5] = 3
$arr[0,6] = 0
$arr[0,7] = "H"
$arr[0,8] = 2
$arr[0,9] = 3.4
$arr[0,10] = 3.1
$arr[0,11] = 2.1
$arr[0,12] = 3.4
$arr[0,13] = 2.9
$arr[0,14] = -0.25
$arr[0,15] = 1.9
$arr[0,16] = 1.9
$arr[0,17] = 2.5
$arr[0,18] = 1.925
$arr[0,19] = 1.875
$arr[0,20] = 1.1
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.8999999999999999
$arr[0,24] = -1
$arr[0,25] = 0.925
$arr[0,26] = -1
$ws.Range("B89:AB89").Value2 = $arr

# Row 90
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845314
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45293.5
$arr[0,3] = "Airdrieonians"
$arr[0,4] = "Inverness CT"
$arr[0,5] = 2
$arr[0,6] = 0
$arr[0,7] = "H"
$arr[0,8] = 2.2
$arr[0,9] = 3.5
$arr[0,10] = 2.7
$arr[0,11] = 2.3
$arr[0,12] = 3.25
$arr[0,13] = 2.7
$arr[0,14] = 0
$arr[0,15] = 1.75
$arr[0,16] = 2.05
$arr[0,17] = 2.25
$arr[0,18] = 1.925
$arr[0,19] = 1.925
$arr[0,20] = 1.3
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.75
$arr[0,24] = -1
$arr[0,25] = -0.5
$arr[0,26] = 0.4625
$ws.Range("B90:AB90").Value2 = $arr

# Row 91
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845311
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45293.5
$arr[0,3] = "Dunfermline"
$arr[0,4] = "Raith"
$arr[0,5] = 1
$arr[0,6] = 2
$arr[0,7] = "A"
$arr[0,8] = 2.6
$arr[0,9] = 3.5
$arr[0,10] = 2.25
$arr[0,11] = 3.2
$arr[0,12] = 3.75
$arr[0,13] = 2
$arr[0,14] = 0.25
$arr[0,15] = 2
$arr[0,16] = 1.8
$arr[0,17] = 2.75
$arr[0,18] = 1.925
$arr[0,19] = 1.875
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 1
$arr[0,23] = -1
$arr[0,24] = 0.8
$arr[0,25] = 0.4625
$arr[0,26] = -0.5
$ws.Range("B91:AB91").Value2 = $arr

# Row 92
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845313
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45293.5
$arr[0,3] = "Partick"
$arr[0,4] = "Queens Park"
$arr[0,5] = 3
$arr[0,6] = 2
$arr[0,7] = "H"
$arr[0,8] = 1.666
$arr[0,9] = 3.75
$arr[0,10] = 4
$arr[0,11] = 1.4
$arr[0,12] = 4.75
$arr[0,13] = 5.75
$arr[0,14] = -1.25
$arr[0,15] = 1.875
$arr[0,16] = 1.975
$arr[0,17] = 3.25
$arr[0,18] = 1.95
$arr[0,19] = 1.9
$arr[0,20] = 0.3999999999999999
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = -0.5
$arr[0,24] = 0.4875
$arr[0,25] = 0.95
$arr[0,26] = -1
$ws.Range("B92:AB92").Value2 = $arr

# Row 115
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845336
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45339.5
$arr[0,3] = "Dunfermline"
$arr[0,4] = "Arbroath"
$arr[0,5] = 1
$arr[0,6] = 1
$arr[0,7] = "D"
$arr[0,8] = 1.5
$arr[0,9] = 4
$arr[0,10] = 5.5
$arr[0,11] = 1.666
$arr[0,12] = 3.6
$arr[0,13] = 4.5
$arr[0,14] = -0.75
$arr[0,15] = 1.95
$arr[0,16] = 1.85
$arr[0,17] = 2.5
$arr[0,18] = 1.9
$arr[0,19] = 1.9
$arr[0,20] = -1
$arr[0,21] = 2.6
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.8500000000000001
$arr[0,25] = -1
$arr[0,26] = 0.8999999999999999
$ws.Range("B115:AB115").Value2 = $arr

# Row 116
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845337
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45339.5
$arr[0,3] = "Inverness CT"
$arr[0,4] = "Partick"
$arr[0,5] = 3
$arr[0,6] = 3
$arr[0,7] = "D"
$arr[0,8] = 2.6
$arr[0,9] = 3.4
$arr[0,10] = 2.4
$arr[0,11] = 2.55
$arr[0,12] = 3.3
$arr[0,13] = 2.5
$arr[0,14] = 0
$arr[0,15] = 1.975
$arr[0,16] = 1.875
$arr[0,17] = 2.5
$arr[0,18] = 2
$arr[0,19] = 1.85
$arr[0,20] = -1
$arr[0,21] = 2.3
$arr[0,22] = -1
$arr[0,23] = 0
$arr[0,24] = 0
$arr[0,25] = 1
$arr[0,26] = -1
$ws.Range("B116:AB116").Value2 = $arr

# Row 117
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6845338
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45339.5
$arr[0,3] = "Queens Park"
$arr[0,4] = "Morton"
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = "D"
$arr[0,8] = 2.4
$arr[0,9] = 3.6
$arr[0,10] = 2.4
$arr[0,11] = 3
$arr[0,12] = 3.6
$arr[0,13] = 2.05
$arr[0,14] = 0.25
$arr[0,15] = 1.975
$arr[0,16] = 1.875
$arr[0,17] = 2.5
$arr[0,18] = 1.925
$arr[0,19] = 1.925
$arr[0,20] = -1
$arr[0,21] = 2.6
$arr[0,22] = -1
$arr[0,23] = 0.4875
$arr[0,24] = -0.5
$arr[0,25] = -1
$arr[0,26] = 0.925
$ws.Range("B117:AB117").Value2 = $arr

# Row 143
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6957817
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45367.5
$arr[0,3] = "Partick"
$arr[0,4] = "Morton"
$arr[0,5] = 2
$arr[0,6] = 1
$arr[0,7] = "H"
$arr[0,8] = 1.8
$arr[0,9] = 3.6
$arr[0,10] = 3.8
$arr[0,11] = 2.3
$arr[0,12] = 3.1
$arr[0,13] = 3
$arr[0,14] = -0.25
$arr[0,15] = 2.025
$arr[0,16] = 1.775
$arr[0,17] = 2.25
$arr[0,18] = 1.9
$arr[0,19] = 1.9
$arr[0,20] = 1.3
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 1.025
$arr[0,24] = -1
$arr[0,25] = 0.8999999999999999
$arr[0,26] = -1
$ws.Range("B143:AB143").Value2 = $arr

# Row 144
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6975416
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45367.5
$arr[0,3] = "Inverness CT"
$arr[0,4] = "Ayr"
$arr[0,5] = 1
$arr[0,6] = 2
$arr[0,7] = "A"
$arr[0,8] = 2.1
$arr[0,9] = 3.4
$arr[0,10] = 3.1
$arr[0,11] = 2.15
$arr[0,12] = 3.1
$arr[0,13] = 3.2
$arr[0,14] = -0.25
$arr[0,15] = 1.9
$arr[0,16] = 1.9
$arr[0,17] = 2.25
$arr[0,18] = 1.9
$arr[0,19] = 1.9
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 2.2
$arr[0,23] = -1
$arr[0,24] = 0.8999999999999999
$arr[0,25] = 0.8999999999999999
$arr[0,26] = -1
$ws.Range("B144:AB144").Value2 = $arr

# Row 145
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6994673
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45374.5
$arr[0,3] = "Arbroath"
$arr[0,4] = "Partick"
$arr[0,5] = 0
$arr[0,6] = 1
$arr[0,7] = "A"
$arr[0,8] = 4.2
$arr[0,9] = 4
$arr[0,10] = 1.6
$arr[0,11] = 5
$arr[0,12] = 4.333
$arr[0,13] = 1.45
$arr[0,14] = 1.25
$arr[0,15] = 1.825
$arr[0,16] = 1.975
$arr[0,17] = 3
$arr[0,18] = 1.8
$arr[0,19] = 2
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 0.45
$arr[0,23] = 0.4125
$arr[0,24] = -0.5
$arr[0,25] = -1
$arr[0,26] = 1
$ws.Range("B145:AB145").Value2 = $arr

# Row 148
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 6975419
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45374.5
$arr[0,3] = "Morton"
$arr[0,4] = "Dunfermline"
$arr[0,5] = 0
$arr[0,6] = 1
$arr[0,7] = "A"
$arr[0,8] = 1.8
$arr[0,9] = 3.5
$arr[0,10] = 3.5
$arr[0,11] = 2.1
$arr[0,12] = 3.25
$arr[0,13] = 3
$arr[0,14] = -0.25
$arr[0,15] = 1.95
$arr[0,16] = 1.85
$arr[0,17] = 2.25
$arr[0,18] = 1.95
$arr[0,19] = 1.85
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 2
$arr[0,23] = -1
$arr[0,24] = 0.8500000000000001
$arr[0,25] = -1
$arr[0,26] = 0.8500000000000001
$ws.Range("B148:AB148").Value2 = $arr

# Row 150
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 7024072
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45381.5
$arr[0,3] = "Partick"
$arr[0,4] = "Inverness CT"
$arr[0,5] = 1
$arr[0,6] = 0
$arr[0,7] = "H"
$arr[0,8] = 1.727
$arr[0,9] = 3.5
$arr[0,10] = 4.333
$arr[0,11] = 1.8
$arr[0,12] = 3.4
$arr[0,13] = 4.2
$arr[0,14] = -0.5
$arr[0,15] = 1.85
$arr[0,16] = 1.95
$arr[0,17] = 2.5
$arr[0,18] = 1.975
$arr[0,19] = 1.825
$arr[0,20] = 0.8
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.8500000000000001
$arr[0,24] = -1
$arr[0,25] = -1
$arr[0,26] = 0.825
$ws.Range("B150:AB150").Value2 = $arr

# Row 151
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 7024061
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45381.5
$arr[0,3] = "Morton"
$arr[0,4] = "Queens Park"
$arr[0,5] = 2
$arr[0,6] = 0
$arr[0,7] = "H"
$arr[0,8] = 1.909
$arr[0,9] = 3.4
$arr[0,10] = 3.6
$arr[0,11] = 1.95
$arr[0,12] = 3.4
$arr[0,13] = 3.5
$arr[0,14] = -0.5
$arr[0,15] = 1.975
$arr[0,16] = 1.825
$arr[0,17] = 2.5
$arr[0,18] = 1.975
$arr[0,19] = 1.825
$arr[0,20] = 0.95
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.9750000000000001
$arr[0,24] = -1
$arr[0,25] = -1
$arr[0,26] = 0.825
$ws.Range("B151:AB151").Value2 = $arr

# Row 152
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 7020832
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45381.5
$arr[0,3] = "Arbroath"
$arr[0,4] = "Dunfermline"
$arr[0,5] = 2
$arr[0,6] = 3
$arr[0,7] = "A"
$arr[0,8] = 3.75
$arr[0,9] = 3.5
$arr[0,10] = 1.833
$arr[0,11] = 4.333
$arr[0,12] = 3.6
$arr[0,13] = 1.7
$arr[0,14] = 0.75
$arr[0,15] = 1.9
$arr[0,16] = 1.9
$arr[0,17] = 2.5
$arr[0,18] = 1.9
$arr[0,19] = 1.9
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 0.7
$arr[0,23] = -0.5
$arr[0,24] = 0.45
$arr[0,25] = 0.8999999999999999
$arr[0,26] = -1
$ws.Range("B152:AB152").Value2 = $arr

# Row 153
$arr = New-Object "object[,]" 1,27
$arr[0,0] = 7020833
$arr[0,1] = "Scotland Championship"
$arr[0,2] = 45381.5
$arr[0,3] = "Dundee Utd"
$arr[0,4] = "Raith"
$arr[0,5] = 2
$arr[0,6] = 0
$arr[0,7] = "H"
$arr[0,8] = 1.666
$arr[0,9] = 3.75
$arr[0,10] = 4.333
$arr[0,11] = 1.75
$arr[0,12] = 3.6
$arr[0,13] = 4
$arr[0,14] = -0.5
$arr[0,15] = 1.825
$arr[0,16] = 1.975
$arr[0,17] = 2.5
$arr[0,18] = 1.975
$arr[0,19] = 1.825
$arr[0,20] = 0.75
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.825
$arr[0,24] = -1
$arr[0,25] = -1
$arr[0,26] = 0.825
$ws.Range("B153:AB153").Value2 = $arr

# Rows 173-176: direct cell edits (upcoming fixtures updated)
# Row 173
$ws.Range("D173").Value2 = 45415.65625
$ws.Range("E173").Value2 = "Raith"
$ws.Range("F173").Value2 = "Arbroath"
$ws.Range("J173").Value2 = 1.363
$ws.Range("K173").Value2 = 4.75
$ws.Range("L173").Value2 = 6.5
$ws.Range("M173").Value2 = 1.333
$ws.Range("N173").Value2 = 4.5
$ws.Range("O173").Value2 = 8
$ws.Range("P173").Value2 = -1.5
$ws.Range("Q173").Value2 = 2
$ws.Range("R173").Value2 = 1.85
$ws.Range("S173").Value2 = 3
$ws.Range("T173").Value2 = 2.025
$ws.Range("U173").Value2 = 1.825

# Row 174
$ws.Range("D174").Value2 = 45415.65625
$ws.Range("E174").Value2 = "Ayr"
$ws.Range("F174").Value2 = "Dunfermline"
$ws.Range("J174").Value2 = 2.5
$ws.Range("K174").Value2 = 3.4
$ws.Range("L174").Value2 = 2.5
$ws.Range("M174").Value2 = 2.375
$ws.Range("N174").Value2 = 3.4
$ws.Range("O174").Value2 = 2.625
$ws.Range("P174").Value2 = 0
$ws.Range("Q174").Value2 = 1.825
$ws.Range("R174").Value2 = 2.025
$ws.Range("S174").Value2 = 2.5
$ws.Range("T174").Value2 = 1.95
$ws.Range("U174").Value2 = 1.9

# Row 175
$ws.Range("D175").Value2 = 45415.65625
$ws.Range("E175").Value2 = "Dundee Utd"
$ws.Range("F175").Value2 = "Partick"
$ws.Range("J175").Value2 = 1.444
$ws.Range("K175").Value2 = 4.2
$ws.Range("L175").Value2 = 6
$ws.Range("M175").Value2 = 1.533
$ws.Range("N175").Value2 = 4
$ws.Range("O175").Value2 = 5.25
$ws.Range("P175").Value2 = -1
$ws.Range("Q175").Value2 = 1.875
$ws.Range("R175").Value2 = 1.975
$ws.Range("S175").Value2 = 3
$ws.Range("T175").Value2 = 2
$ws.Range("U175").Value2 = 1.85

# Row 176
$ws.Range("D176").Value2 = 45415.65625
$ws.Range("E176").Value2 = "Inverness CT"
$ws.Range("F176").Value2 = "Morton"
$ws.Range("J176").Value2 = 2.1
$ws.Range("K176").Value2 = 3.4
$ws.Range("L176").Value2 = 3.1
$ws.Range("M176").Value2 = 2.05
$ws.Range("N176").Value2 = 3.4
$ws.Range("O176").Value2 = 3.2
$ws.Range("P176").Value2 = -0.25
$ws.Range("Q176").Value2 = 1.85
$ws.Range("R176").Value2 = 2
$ws.Range("S176").Value2 = 2.25
$ws.Range("T176").Value2 = 1.75
$ws.Range("U176").Value2 = 2.125

# Update match-id values (B173:B176)
$ws.Range("B173").Value2 = 7222798
$ws.Range("B174").Value2 = 7248178
$ws.Range("B175").Value2 = 7251789
$ws.Range("B176").Value2 = 7251790